$d = $word.ActiveDocument

# --- Edit 2: "sprinkle again with " -> "sprinkle " + "once " + "again with " ---
# The new word "once " must sit in its own run that has no explicit color
# (just rtl=0), matching the formatting used elsewhere in the document for
# "automatic" colored text (e.g. the lone "M" run near the top of the
# document, in "<head>Molding"). A plain Find/Replace would just extend the
# existing colored run it lands in, so instead we borrow the formatting
# from that plain "M" run via FormattedText, drop it in at the insertion
# point, then set its text. This must run before Edit 1 below, since Edit 1
# rewrites that very "M" run and Word's Range/FormattedText objects stay
# live/bound to that text, so capturing after Edit 1 would pick up the
# wrong (already-changed) formatting.

$plainTextIdx = $d.Content.Text.IndexOf("Molding")
$plainFormatSrc = $d.Range($plainTextIdx, $plainTextIdx + 1)
$plainFormat = $plainFormatSrc.FormattedText

$r = $d.Content
$r.Find.Execute("sprinkle ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)   # wdCollapseEnd

$insStart = $r.Start
$r.InsertAfter("X")
$insRange = $d.Range($insStart, $insStart + 1)
$insRange.FormattedText = $plainFormat
$insRange.Text = "once "

# --- Edit 1: "<head>" + "M" + "olding " -> "<head>For m" + "olding " ---
# Replace "<head>M" with "<head>For m"; Word preserves the formatting of the
# run in which the match starts (the "<head>" run, Courier New / 7f6000),
# merging the result into a single run and leaving "olding " untouched.
$d.Content.Find.Execute("<head>M", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<head>For m", 2) | Out-Null
